$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Atelier")

# Column A: shift all dates from 2015xxxx -> 2017xxxx (rows 3-63)
for ($r = 3; $r -le 63; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $cell.Value2 = $cell.Value2 + 20000
}

# Column E: updated values per row (rows with changes only)
$eValues = @{
    3  = 12
    5  = 6
    6  = 15
    8  = 8
    9  = 14
    10 = 13
    11 = 11
    12 = 17
    13 = 7
    14 = 8
    15 = 15
    16 = 15
    17 = 14
    18 = 18
    19 = 12
    21 = 8
    22 = 16
    23 = 18
    24 = 6
    25 = 9
    26 = 5
    27 = 15
    28 = 17
    29 = 8
    30 = 18
    31 = 11
    32 = 16
    33 = 9
    34 = 7
    35 = 6
    36 = 8
    37 = 17
    38 = 6
    39 = 13
    40 = 20
    41 = 15
    42 = 15
    43 = 8
    44 = 12
    45 = 17
    46 = 14
    47 = 17
    48 = 16
    49 = 19
    50 = 13
    51 = 6
    52 = 5
    53 = 9
    54 = 5
    55 = 18
    56 = 7
    57 = 10
    58 = 7
    59 = 15
    60 = 5
    61 = 8
    62 = 6
    63 = 7
}

foreach ($row in $eValues.Keys) {
    $ws.Cells.Item($row, 5).Value2 = $eValues[$row]
}
